{"js": "// Replace the date and each \"AxB=C\" answer cell with its updated value.\n// Every \"old\" search string is unique in the document, so a single\n// search+replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-06-13 Thursday\", \"2024-06-14 Friday\"],\n  [\"62\u00d759=3658\", \"37\u00d736=1332\"],\n  [\"84\u00d740=3360\", \"27\u00d741=1107\"],\n  [\"90\u00d789=8010\", \"86\u00d795=8170\"],\n  [\"54\u00d774=3996\", \"96\u00d719=1824\"],\n  [\"33\u00d731=1023\", \"18\u00d765=1170\"],\n  [\"25\u00d798=2450\", \"33\u00d750=1650\"],\n  [\"75\u00d711=825\", \"78\u00d725=1950\"],\n  [\"81\u00d762=5022\", \"29\u00d760=1740\"],\n  [\"58\u00d770=4060\", \"78\u00d788=6864\"],\n  [\"11\u00d723=253\", \"18\u00d765=1170\"],\n  [\"42\u00d771=2982\", \"77\u00d759=4543\"],\n  [\"49\u00d727=1323\", \"78\u00d736=2808\"],\n  [\"29\u00d779=2291\", \"84\u00d760=5040\"],\n  [\"24\u00d716=384\", \"74\u00d721=1554\"],\n  [\"87\u00d794=8178\", \"91\u00d733=3003\"],\n  [\"56\u00d782=4592\", \"67\u00d789=5963\"],\n  [\"54\u00d770=3780\", \"38\u00d723=874\"],\n  [\"85\u00d799=8415\", \"80\u00d728=2240\"],\n  [\"97\u00d775=7275\", \"19\u00d795=1805\"],\n  [\"76\u00d713=988\", \"91\u00d760=5460\"],\n  [\"17\u00d724=408\", \"61\u00d719=1159\"],\n  [\"14\u00d794=1316\", \"72\u00d762=4464\"],\n  [\"60\u00d745=2700\", \"23\u00d771=1633\"],\n  [\"91\u00d774=6734\", \"54\u00d790=4860\"],\n  [\"37\u00d734=1258\", \"77\u00d714=1078\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each tuple is (search text, replacement text). All search strings are\n# unique within the document, so a single Find/Replace per pair is enough.\n$pairs = @(\n    @(\"2024-06-13 Thursday\", \"2024-06-14 Friday\"),\n    @(\"62\u00d759=3658\", \"37\u00d736=1332\"),\n    @(\"84\u00d740=3360\", \"27\u00d741=1107\"),\n    @(\"90\u00d789=8010\", \"86\u00d795=8170\"),\n    @(\"54\u00d774=3996\", \"96\u00d719=1824\"),\n    @(\"33\u00d731=1023\", \"18\u00d765=1170\"),\n    @(\"25\u00d798=2450\", \"33\u00d750=1650\"),\n    @(\"75\u00d711=825\", \"78\u00d725=1950\"),\n    @(\"81\u00d762=5022\", \"29\u00d760=1740\"),\n    @(\"58\u00d770=4060\", \"78\u00d788=6864\"),\n    @(\"11\u00d723=253\", \"18\u00d765=1170\"),\n    @(\"42\u00d771=2982\", \"77\u00d759=4543\"),\n    @(\"49\u00d727=1323\", \"78\u00d736=2808\"),\n    @(\"29\u00d779=2291\", \"84\u00d760=5040\"),\n    @(\"24\u00d716=384\", \"74\u00d721=1554\"),\n    @(\"87\u00d794=8178\", \"91\u00d733=3003\"),\n    @(\"56\u00d782=4592\", \"67\u00d789=5963\"),\n    @(\"54\u00d770=3780\", \"38\u00d723=874\"),\n    @(\"85\u00d799=8415\", \"80\u00d728=2240\"),\n    @(\"97\u00d775=7275\", \"19\u00d795=1805\"),\n    @(\"76\u00d713=988\", \"91\u00d760=5460\"),\n    @(\"17\u00d724=408\", \"61\u00d719=1159\"),\n    @(\"14\u00d794=1316\", \"72\u00d762=4464\"),\n    @(\"60\u00d745=2700\", \"23\u00d771=1633\"),\n    @(\"91\u00d774=6734\", \"54\u00d790=4860\"),\n    @(\"37\u00d734=1258\", \"77\u00d714=1078\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}\n"}
